# Update the dSF column (F) values for several rows.
# These correspond to a "repull data" of the dSF figures that no longer
# mirror the dS0 column (E) as they did before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -6
$ws.Range("F4").Value  = -6
$ws.Range("F5").Value  = -2
$ws.Range("F7").Value  = -5
$ws.Range("F8").Value  = -4
$ws.Range("F18").Value = 2
$ws.Range("F19").Value = -2
